$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the data rows ---
# Row 2 ("sgdgf" -> "chill")
$ws.Range("B2").Value = "chill"
# Row 3 ("Rohan" -> "shit", price 535.6 -> 300)
$ws.Range("B3").Value = "shit"
$ws.Range("C3").Value = 300

# Remove the trailing rows (old rows 4-6: chill/Sohan/fooking) - no longer needed
$ws.Rows("4:6").Delete() | Out-Null

# --- Strip the bold "Heading 1" styling from the header row back to Normal ---
$ws.Range("A1:C1").Style = "Normal"
# Let the row re-measure its height now that the big bold font is gone
$ws.Rows("1:1").AutoFit() | Out-Null
# Remove the now-unused named cell style from the workbook's style gallery
try { $wb.Styles("Heading 1").Delete() | Out-Null } catch {}

# --- Match the saved selection / active cell ---
$ws.Range("C1").Select() | Out-Null
